$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column S header (year 2022) matching style of R4
$ws.Range("R4").Copy()
$ws.Range("S4").PasteSpecial(-4122)
$ws.Range("S4").Value = 2022

# Update existing data values in row 5
$ws.Range("P5").Value = 20.5
$ws.Range("Q5").Value = 20.5
$ws.Range("R5").Value = 17.9

# New value for S5 matching style of R5
$ws.Range("R5").Copy()
$ws.Range("S5").PasteSpecial(-4122)
$ws.Range("S5").Value = 13.5

# Update selection to reflect new active cell / selection range
$ws.Range("S7:S8").Select()
